$d = $word.ActiveDocument

# Helper: insert a WordprocessingML body fragment at a given Range by
# wrapping it in a minimal OOXML package, as InsertXML expects.
function Insert-WordXmlAt($range, [string]$innerBody) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $result = $range.InsertXML($pkg)
}

# The trailing "_GoBack" bookmark currently sits right after "Arreglar css";
# it needs to end up at the very end of the document once the new content
# has been appended below, so drop it here and recreate it there later.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Split the "css" word out of the "Arreglar css" bullet into its own run
# and wrap it with spell-check proofErr markers (mirrors the diff, which
# turns "Arreglar css" into "Arreglar " + proofErr-wrapped "css").
$searchRange = $d.Content
$searchRange.Find.Execute("css") | Out-Null
# Rebuild a plain Range from the hit so InsertXML replaces just that text
# instead of appending after it.
$cssRange = $d.Range($searchRange.Start, $searchRange.End)
Insert-WordXmlAt $cssRange '<w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>css</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

# Append the new paragraphs at the end of the document: a blank line, the
# "Clientes" heading, its two bullet items, a blank heading-styled line,
# and a final blank paragraph that will hold the relocated bookmark.
$endRange = $d.Range($d.Content.End, $d.Content.End)
$tail = '<w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr></w:p>' +
    '<w:p><w:pPr><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>Clientes</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Fecha del día de hoy cuando al crear un cliente</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Colocar campos formulario.</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr></w:p>' +
    '<w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Insert-WordXmlAt $endRange $tail
